$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on all target cells first so numeric-looking strings
# (e.g. "174.10", "0.998") are stored as text, matching the source inlineStr type,
# rather than being auto-converted to numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("B18").NumberFormat = "@"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("B19").NumberFormat = "@"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"

$ws.Range("D2").Value = '69.003.12'
$ws.Range("E2").Value = '  -3.95%  '
$ws.Range("D3").Value = '3.514.13'
$ws.Range("E3").Value = '  -3.71%  '
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").Value = '576.92'
$ws.Range("E5").Value = '  -2.44%  '
$ws.Range("D6").Value = '174.10'
$ws.Range("E6").Value = '  -4.05%  '
$ws.Range("E7").Value = '  -0.21%  '
$ws.Range("D8").Value = '3.504.13'
$ws.Range("E8").Value = '  -3.78%  '
$ws.Range("E9").Value = '  -0.03%  '
$ws.Range("E10").Value = '  -7.02%  '
$ws.Range("D11").Value = '6.65'
$ws.Range("E11").Value = '  +8.80%  '
$ws.Range("D12").Value = '0.605'
$ws.Range("E12").Value = '  -0.86%  '
$ws.Range("D13").Value = '47.25'
$ws.Range("E13").Value = '  -5.44%  '
$ws.Range("E14").Value = '  -4.12%  '
$ws.Range("D15").Value = '687.42'
$ws.Range("E15").Value = '  +0.84%  '
$ws.Range("D16").Value = '8.91'
$ws.Range("E16").Value = '  -1.66%  '
$ws.Range("D17").Value = '4.089.27'
$ws.Range("E17").Value = '  -3.33%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '69.138.25'
$ws.Range("E18").Value = '  -3.83%  '
$ws.Range("B19").Value = 'WrappedEther'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D19").Value = '3.518.55'
$ws.Range("E19").Value = '  -3.64%  '
$ws.Range("E20").Value = '  -1.42%  '
$ws.Range("D21").Value = '17.51'
$ws.Range("E21").Value = '  -4.14%  '
$ws.Range("D22").Value = '11.21'
$ws.Range("E22").Value = '  -3.85%  '
$ws.Range("E23").Value = '  -3.82%  '
$ws.Range("D24").Value = '16.53'
$ws.Range("E24").Value = '  -7.92%  '
$ws.Range("D25").Value = '97.57'
$ws.Range("E25").Value = '  -5.62%  '
$ws.Range("E26").Value = '  -4.87%  '
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("D28").Value = '2.66'
$ws.Range("E28").Value = '  -6.40%  '
$ws.Range("D29").Value = '9.39'
$ws.Range("E29").Value = '  -7.67%  '
$ws.Range("D30").Value = '33.10'
$ws.Range("E30").Value = '  -6.47%  '
$ws.Range("D31").Value = '8.83'
$ws.Range("E31").Value = '  -4.18%  '
$ws.Range("D32").Value = '3.17'
$ws.Range("E32").Value = '  -8.70%  '
$ws.Range("D33").Value = '7.29'
$ws.Range("E33").Value = '  -0.89%  '
$ws.Range("E34").Value = '  -5.78%  '
$ws.Range("D35").Value = '562.86'
$ws.Range("E35").Value = '  -2.94%  '
$ws.Range("D36").Value = '3.65'
$ws.Range("E36").Value = '  -13.09%  '
$ws.Range("D37").Value = '10.87'
$ws.Range("E37").Value = '  -4.24%  '
$ws.Range("D38").Value = '0.106'
$ws.Range("E38").Value = '  -3.38%  '
$ws.Range("D39").Value = '57.28'
$ws.Range("E39").Value = '  -3.64%  '
$ws.Range("D40").Value = '0.998'
$ws.Range("E40").Value = '  +0.02%  '
$ws.Range("E41").Value = '  -3.83%  '
$ws.Range("D42").Value = '0.0443'
$ws.Range("E42").Value = '  -5.68%  '
$ws.Range("D43").Value = '3.461.00'
$ws.Range("E43").Value = '  -7.20%  '
$ws.Range("D44").Value = '0.337'
$ws.Range("E44").Value = '  -2.87%  '
$ws.Range("D45").Value = '33.45'
$ws.Range("E45").Value = '  -6.48%  '
$ws.Range("E46").Value = '  -8.11%  '
$ws.Range("D47").Value = '2.92'
$ws.Range("E47").Value = '  +3.23%  '
$ws.Range("E48").Value = '  -8.11%  '
$ws.Range("E49").Value = '  -0.20%  '
$ws.Range("D50").Value = '134.32'
$ws.Range("E50").Value = '  +1.36%  '
$ws.Range("E51").Value = '  -0.36%  '

# Clear the temporary text-number-format styling so saved cells have no
# style index, matching the original (unstyled) data cells.
$ws.Range("D2").ClearFormats()
$ws.Range("E2").ClearFormats()
$ws.Range("D3").ClearFormats()
$ws.Range("E3").ClearFormats()
$ws.Range("E4").ClearFormats()
$ws.Range("D5").ClearFormats()
$ws.Range("E5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("E6").ClearFormats()
$ws.Range("E7").ClearFormats()
$ws.Range("D8").ClearFormats()
$ws.Range("E8").ClearFormats()
$ws.Range("E9").ClearFormats()
$ws.Range("E10").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("E11").ClearFormats()
$ws.Range("D12").ClearFormats()
$ws.Range("E12").ClearFormats()
$ws.Range("D13").ClearFormats()
$ws.Range("E13").ClearFormats()
$ws.Range("E14").ClearFormats()
$ws.Range("D15").ClearFormats()
$ws.Range("E15").ClearFormats()
$ws.Range("D16").ClearFormats()
$ws.Range("E16").ClearFormats()
$ws.Range("D17").ClearFormats()
$ws.Range("E17").ClearFormats()
$ws.Range("B18").ClearFormats()
$ws.Range("C18").ClearFormats()
$ws.Range("D18").ClearFormats()
$ws.Range("E18").ClearFormats()
$ws.Range("B19").ClearFormats()
$ws.Range("C19").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("E19").ClearFormats()
$ws.Range("E20").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("E21").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("E22").ClearFormats()
$ws.Range("E23").ClearFormats()
$ws.Range("D24").ClearFormats()
$ws.Range("E24").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("E25").ClearFormats()
$ws.Range("E26").ClearFormats()
$ws.Range("E27").ClearFormats()
$ws.Range("D28").ClearFormats()
$ws.Range("E28").ClearFormats()
$ws.Range("D29").ClearFormats()
$ws.Range("E29").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("E30").ClearFormats()
$ws.Range("D31").ClearFormats()
$ws.Range("E31").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("E32").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("E33").ClearFormats()
$ws.Range("E34").ClearFormats()
$ws.Range("D35").ClearFormats()
$ws.Range("E35").ClearFormats()
$ws.Range("D36").ClearFormats()
$ws.Range("E36").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("E37").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("E38").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("E39").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("E40").ClearFormats()
$ws.Range("E41").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("E42").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("E43").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("E44").ClearFormats()
$ws.Range("D45").ClearFormats()
$ws.Range("E45").ClearFormats()
$ws.Range("E46").ClearFormats()
$ws.Range("D47").ClearFormats()
$ws.Range("E47").ClearFormats()
$ws.Range("E48").ClearFormats()
$ws.Range("E49").ClearFormats()
$ws.Range("D50").ClearFormats()
$ws.Range("E50").ClearFormats()
$ws.Range("E51").ClearFormats()
